$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "false start" data rows (old rows 2 and 3), shifting the
# remaining rows up so the data that used to be in rows 4 and 5 becomes
# the new rows 2 and 3.
$rows = $ws.Rows("2:3")
$rows.Select() | Out-Null
$rows.Delete() | Out-Null
